$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix tiny floating-point precision on the existing last row's timestamp (row 30)
$ws.Range("A30").Value = 44343.79602429167

# Append new data row 31 (new job-numbers snapshot)
$ws.Range("A31").Value = 44344.80417441898
$ws.Range("B31").Value = 74878
$ws.Range("C31").Value = 62946
$ws.Range("D31").Value = 3222
$ws.Range("E31").Value = 2056
$ws.Range("F31").Value = 1430
$ws.Range("G31").Value = 19490
$ws.Range("H31").Value = 1360
$ws.Range("I31").Value = 817
$ws.Range("J31").Value = 197

$wb.Save()
